# Week 15 simulations update
$wb = $excel.ActiveWorkbook

# --- K sheet: kicker M.Ammendola replaced by E.Pineiro (same row/stats) ---
$wsK = $wb.Worksheets.Item("K")
$wsK.Range("A2").Value = "E.Pineiro"

# --- WR sheet: new player D.Montgomery added as row 10 (all stats 0) ---
$wsWR = $wb.Worksheets.Item("WR")
$wsWR.Cells.Item(10, 1).Value = "D.Montgomery"
$wsWR.Cells.Item(10, 2).Value = 0
$wsWR.Cells.Item(10, 3).Value = 0
$wsWR.Cells.Item(10, 4).Value = 0
$wsWR.Cells.Item(10, 5).Value = 0
$wsWR.Cells.Item(10, 6).Value = 0
$wsWR.Cells.Item(10, 7).Value = 0
$wsWR.Cells.Item(10, 8).Value = 0
$wsWR.Cells.Item(10, 9).Value = 0
$wsWR.Cells.Item(10, 10).Value = 0

# --- Update selections left behind on individual sheets ---
# K sheet selection moves to C4 (sheet is activated momentarily to set this,
# matching real Excel behaviour where Select() requires the sheet be active)
$wsK.Activate()
$wsK.Range("C4").Select()

# --- Final active sheet/tab is WR, with its own selection at O25 ---
$wsWR.Activate()
$wsWR.Range("O25").Select()
